# DSA Linked List Practice — add 7 new rows (69-75) of linked-list questions,
# matching the format/style used by the existing table (rows 3-68).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 71: 234. Palindrome Linked List -------------------------------
$ws.Range("A71").Value = 69
$ws.Range("B71").Value = "234. Palindrome Linked List"
$ws.Range("C71").Value = "LeetCode"
$ws.Range("D71").Value = "Easy"
$ws.Range("E71").Value = "https://leetcode.com/problems/palindrome-linked-list/description/"
$ws.Range("F71").Value = "|"
$c71 = $ws.Cells.Item(71, 5)
$ws.Hyperlinks.Add($c71, "https://leetcode.com/problems/palindrome-linked-list/description/") | Out-Null

# --- Row 72: 328. Odd Even Linked List ----------------------------------
$ws.Range("A72").Value = 70
$ws.Range("B72").Value = "328. Odd Even Linked List"
$ws.Range("C72").Value = "LeetCode"
$ws.Range("D72").Value = "Medium"
$ws.Range("E72").Value = "https://leetcode.com/problems/odd-even-linked-list/description/"
$ws.Range("F72").Value = "|"
$c72 = $ws.Cells.Item(72, 5)
$ws.Hyperlinks.Add($c72, "https://leetcode.com/problems/odd-even-linked-list/description/") | Out-Null

# --- Row 74: 2095. Delete the Middle Node of a Linked List --------------
$ws.Range("A74").Value = 72
$ws.Range("B74").Value = "2095. Delete the Middle Node of a Linked List"
$ws.Range("C74").Value = "LeetCode"
$ws.Range("D74").Value = "Medium"
$ws.Range("E74").Value = "https://leetcode.com/problems/delete-the-middle-node-of-a-linked-list/description/"
$ws.Range("F74").Value = "|"
$c74 = $ws.Cells.Item(74, 5)
$ws.Hyperlinks.Add($c74, "https://leetcode.com/problems/delete-the-middle-node-of-a-linked-list/description/") | Out-Null

# --- Row 73: 19. Remove Nth Node From End of List ------------------------
# (URL was entered before the question title for this row.)
$ws.Range("A73").Value = 71
$ws.Range("E73").Value = "https://leetcode.com/problems/remove-nth-node-from-end-of-list/description/"
$ws.Range("B73").Value = "19. Remove Nth Node From End of List"
$ws.Range("C73").Value = "LeetCode"
$ws.Range("D73").Value = "Medium"
$ws.Range("F73").Value = "|"
$c73 = $ws.Cells.Item(73, 5)
$ws.Hyperlinks.Add($c73, "https://leetcode.com/problems/remove-nth-node-from-end-of-list/description/") | Out-Null

# --- Row 75: 148. Sort List ----------------------------------------------
$ws.Range("A75").Value = 73
$ws.Range("B75").Value = "148. Sort List"
$ws.Range("C75").Value = "LeetCode"
$ws.Range("D75").Value = "M"
$ws.Range("E75").Value = "https://leetcode.com/problems/sort-list/"
$ws.Range("F75").Value = "|"
$c75 = $ws.Cells.Item(75, 5)
$ws.Hyperlinks.Add($c75, "https://leetcode.com/problems/sort-list/") | Out-Null

# --- Row 70: Find length of Loop -----------------------------------------
$ws.Range("A70").Value = 68
$ws.Range("B70").Value = "Find length of Loop"
$ws.Range("C70").Value = "Naukri Code 360"
$ws.Range("D70").Value = "Easy"
$ws.Range("E70").Value = "https://www.naukri.com/code360/problems/find-length-of-loop_8160455"
$ws.Range("F70").Value = "|"
$c70 = $ws.Cells.Item(70, 5)
$ws.Hyperlinks.Add($c70, "https://www.naukri.com/code360/problems/find-length-of-loop_8160455") | Out-Null

# --- Row 69: 142. Linked List Cycle II -----------------------------------
$ws.Range("A69").Value = 67
$ws.Range("B69").Value = "142. Linked List Cycle II"
$ws.Range("C69").Value = "LeetCode"
$ws.Range("D69").Value = "Medium"
$ws.Range("E69").Value = "https://leetcode.com/problems/linked-list-cycle-ii/description/"
$ws.Range("F69").Value = "|"
$c69 = $ws.Cells.Item(69, 5)
$ws.Hyperlinks.Add($c69, "https://leetcode.com/problems/linked-list-cycle-ii/description/") | Out-Null

# Re-apply the formatting (borders/fonts/hyperlink style) used by the rest
# of the table onto the new rows — done last, since adding a hyperlink
# resets the cell's style.
$ws.Range("A65:F65").Copy() | Out-Null
$ws.Range("A69:F75").PasteSpecial(-4122) | Out-Null

# Recalculate dependent formulas (H column COUNTIFs + the H1 SUM).
$excel.Calculate()

# Update the frozen-pane scroll position / active selection to mirror
# where the author ended up after adding the new rows.
$win = $excel.ActiveWindow
$win.ScrollRow = 45
$win.ScrollColumn = 9
$ws.Range("D74").Select() | Out-Null
